$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Fill B:D with the same value as A for the "single-letter" rows that only
# had column A populated (A2="A", A18="B", A37="C", ... one row per letter).
$rows = @(2,18,37,60,65,74,83,99,105,115,119,129,139,162,177,179,192,194,199,226,242,250,257,260,262)
foreach ($r in $rows) {
    $val = $ws1.Cells.Item($r, 1).Value2
    $ws1.Cells.Item($r, 2).Value = $val
    $ws1.Cells.Item($r, 3).Value = $val
    $ws1.Cells.Item($r, 4).Value = $val
}

# Turn the header row A1:D1 into an AutoFilter range.
[void]$ws1.Range("A1:D1").AutoFilter()

# Repoint the (now stale) _FilterDatabase name at the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$D`$1"
    }
}

# Add a defined name "A" pointing at Sheet1!$JDY$4.
$wb.Names.Add("A", "=Sheet1!`$JDY`$4")

# Make Sheet1 the active sheet/tab (moves tabSelected from Sheet2 to Sheet1).
$ws1.Activate()

# Update the on-screen selection for Sheet1 and scroll the viewport so row
# 214 is at the top (best-effort -- mirrors topLeftCell="A214" in the XML).
[void]$ws1.Range("A266").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 214
$win.ScrollColumn = 1
